$d = $word.ActiveDocument

# --- Change 1: split the run "{m" into two runs "{" and "m" ---------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $s1 = $rng1.Start
    # isolate the "m" character (2nd char of the match) into its own run by
    # toggling a character attribute on/off; ending value equals the
    # original formatting, so no visible formatting change occurs, but the
    # run gets split in two.
    $mChar = $d.Range($s1 + 1, $s1 + 2)
    $mChar.Font.Bold = $true
    $mChar.Font.Bold = $false
}

# --- Change 2: split the run ")}" into ")" and a new, unformatted "}" -----
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute(")}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $s2 = $rng2.Start
    $e2 = $rng2.End
    # isolate the "}" character (2nd char of the match) into its own run,
    # same trick as above.
    $braceChar = $d.Range($s2 + 1, $e2)
    $braceChar.Font.Bold = $true
    $braceChar.Font.Bold = $false

    # Replace the now-isolated "}" run's content with freshly authored OOXML
    # that carries no run formatting at all (and an explicit
    # xml:space="preserve" on the text), matching a run inserted fresh by a
    # rewriter rather than one inheriting the surrounding character style.
    $braceRange = $d.Range($s2 + 1, $e2)
    $xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $braceRange.InsertXML($xml)
}
